$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: corrected hydrogen demand value for Iron & steel
$ws.Range("B3").Value = 224565.8606566988

# D3: value removed (cell now empty)
$ws.Range("D3").ClearContents()

# C4: corrected methanol value for Chemicals
$ws.Range("C4").Value = 39.2830280879696

# C5: corrected ammonia value for Chemicals
$ws.Range("C5").Value = 2466.356526526049

# Row 7: label changed from "Other" to "Biogas", value updated
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 179.353386267281

# New row 8: "Other" row (re-added below Biogas), matching the styling of A7
$ws.Range("A8").Value = "Other"
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D8").Value = 138.1620075955699
